$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("thermoRxns")

$ws.Range("B3").Value = -184.66
$ws.Range("C3").Value = -183.06
$ws.Range("B4").Value = -17.67
$ws.Range("C4").Value = -16.79
$ws.Range("B6").Value = -13.28
$ws.Range("C6").Value = -6.720000000000001
$ws.Range("B8").Value = -15.53
$ws.Range("C8").Value = -8.449999999999999
$ws.Range("B9").Value = -18.47
$ws.Range("C9").Value = -14.65
$ws.Range("B10").Value = -19.27
$ws.Range("C10").Value = -15.39
$ws.Range("B13").Value = -3.75
$ws.Range("C13").Value = -1.05
$ws.Range("B14").Value = -4.55
$ws.Range("C14").Value = -1.79
$ws.Range("B15").Value = -3.75
$ws.Range("C15").Value = -1.05
$ws.Range("B16").Value = -180.32
$ws.Range("C16").Value = -176.72
$ws.Range("B17").Value = 7.049999999999999
$ws.Range("C17").Value = 13.53
$ws.Range("B18").Value = 7.83
$ws.Range("C18").Value = 14.31
$ws.Range("B19").Value = -2.7
$ws.Range("C19").Value = -1.14
$ws.Range("B20").Value = -4.56
$ws.Range("C20").Value = -2.2
$ws.Range("B21").Value = -5.890000000000001
$ws.Range("C21").Value = -1.97
$ws.Range("B22").Value = -11.95
$ws.Range("C22").Value = -8.030000000000001
$ws.Range("B23").Value = -2.02
$ws.Range("C23").Value = 0.86
$ws.Range("B24").Value = 111.76
$ws.Range("C24").Value = 117.2
$ws.Range("B25").Value = 13.43
$ws.Range("C25").Value = 17.67
$ws.Range("B26").Value = -2.92
$ws.Range("C26").Value = -2.14
$ws.Range("B27").Value = -169.9
$ws.Range("C27").Value = -167.82
$ws.Range("B28").Value = -20.36
$ws.Range("C28").Value = -19.28
$ws.Range("B29").Value = -6
$ws.Range("C29").Value = -4.880000000000001
$ws.Range("B30").Value = 7.44
$ws.Range("C30").Value = 8.26
$ws.Range("B31").Value = -18.92
$ws.Range("C31").Value = -18.02
$ws.Range("B32").Value = 3.88
$ws.Range("C32").Value = 4.62
$ws.Range("B33").Value = 152.57
$ws.Range("C33").Value = 154.21
$ws.Range("B34").Value = -28.05
$ws.Range("C34").Value = -27.19
$ws.Range("B35").Value = -313.2
$ws.Range("C35").Value = -306.4
$ws.Range("B36").Value = -18.1
$ws.Range("C36").Value = -16.5
$ws.Range("B37").Value = -184.66
$ws.Range("C37").Value = -183.06
$ws.Range("B38").Value = -64.84
$ws.Range("C38").Value = -63.72
$ws.Range("B39").Value = 64.45
$ws.Range("C39").Value = 65.67
